# Week 16 logged + season sim from Week 17 onward.
# Updates season-total cells on OFF/DEF/ST/TURNS/PEN, and appends the
# new per-play sample data onto the long space-delimited strings on
# the YDS and ST sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# OFF sheet - season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 192
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = 75
$ws.Range("G2").Value = 59
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 33
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 22

$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 147
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 73
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 52
$ws.Range("J3").Value = 44
$ws.Range("L3").Value = 203
$ws.Range("M3").Value = 132
$ws.Range("Q3").Value = 449

# ---------------------------------------------------------------
# DEF sheet - season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 170
$ws.Range("D2").Value = 8
$ws.Range("F2").Value = 56
$ws.Range("G2").Value = 50
$ws.Range("J2").Value = 25
$ws.Range("O2").Value = 18
$ws.Range("P2").Value = 9

$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 185
$ws.Range("E3").Value = 23
$ws.Range("F3").Value = 111
$ws.Range("G3").Value = 31
$ws.Range("I3").Value = 62
$ws.Range("J3").Value = 52
$ws.Range("L3").Value = 232
$ws.Range("M3").Value = 150
$ws.Range("Q3").Value = 404

# ---------------------------------------------------------------
# ST sheet - season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 86
$ws.Range("D2").Value = 51
$ws.Range("F2").Value = 129
$ws.Range("G2").Value = 124
$ws.Range("J2").Value = 42
$ws.Range("K2").Value = 41
$ws.Range("L2").Value = 33
$ws.Range("M2").Value = 24
$ws.Range("N2").Value = 10
$ws.Range("B3").Value = 42

# ---------------------------------------------------------------
# TURNS sheet - season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 8

# ---------------------------------------------------------------
# PEN sheet - season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 17
$ws.Range("B3").Value = 19

# ---------------------------------------------------------------
# YDS sheet - append per-play samples to the long lists
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Value2 + " 43 3 -1 4 1 -2 2 2 4 -2 7 -2 5 5 11 -1 2 8 18 4 3 -1 1 3 3 -1 3 2 2 1 -3 3"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 9 1 6 23 10 7 3 12 14 7 5 7 7 5 5 37 8 20 39 14"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 6 5 4 2 2 10 57 0 2 1 4 1 0 5 1 7 0 11 10 2"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " 3 10 5 10 5 8 13 7 11 4 3 11 10 24 3 33 4 5 7 16 7 4 7 14 2 14 8"

# ---------------------------------------------------------------
# ST sheet - append per-play samples to the long lists
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = $ws.Range("B4").Value2 + " 66 68 49 68"
$ws.Range("B5").Value = $ws.Range("B5").Value2 + " 21 13 16 31"
$ws.Range("B6").Value = $ws.Range("B6").Value2 + " 16"
$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 39 44 32"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " 0 48 0"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 0"
